# Apply "validated all breadcrumbs data in page data files" edit.
#
# Summary of the change (from the OOXML diff):
#  1. Column E ("-BREAD-" / breadcrumb validity) on many rows is toggled
#     between "VALID" and "UNKNOWN" to reflect a re-validation pass.
#  2. A previously-blank row 250 is filled in with a new page entry:
#     A250 = hyperlink "https://discord.heatlabs.net" (new page), with
#     B250 = "PENDING", C250 = "PENDING", D250 = "UNKNOWN", E250 = "VALID".
#  3. The dependent COUNTIF() summary cells (G3, G7, G11, G15, I15, ...)
#     recalculate automatically from the above - no need to touch them by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose breadcrumb status flips from VALID -> UNKNOWN ---
$rowsToUnknown = @(
    6,7,8,9,10,12,13,15,16,17,19,20,21,24,31,35,39,42,66,67,68,70,76,
    106,115,136,196,199,201,203,220,221,222,225,226,227,232,233,236,
    238,241,244,245
)
foreach ($r in $rowsToUnknown) {
    $ws.Range("E$r").Value = "UNKNOWN"
}

# --- Rows whose breadcrumb status flips from UNKNOWN -> VALID ---
$rowsToValid = @(
    50,51,52,53,54,55,56,57,58,59,61,62,86,87,88,89,99,100,101,102,103,
    104,108,111,112,116,117,119,120,121,122,123,124,125,126,127,128,130,
    131,132,133,134,135,137,138,139,140,141,142,143,144,145,147,149,153,
    159,163,166,167,168,169,170,171,176,183,184,186,212,215
)
foreach ($r in $rowsToValid) {
    $ws.Range("E$r").Value = "VALID"
}

# --- Fill in the new page row (250) with its breadcrumb / index data ---
$ws.Range("B250").Value = "PENDING"
$ws.Range("C250").Value = "PENDING"
$ws.Range("D250").Value = "UNKNOWN"
$ws.Range("E250").Value = "VALID"

# A250 needs both the text/hyperlink AND the same visual style (underlined,
# blue Arial) used by every other page-link cell in column A. Copying the
# format from the row above first means the new style exactly reuses the
# existing shared style instead of Excel fabricating a near-duplicate one.
$ws.Range("A249").Copy($ws.Range("A250"))
$ws.Range("A250").Value = "https://discord.heatlabs.net"
$ws.Hyperlinks.Add($ws.Range("A250"), "https://heatlabs.net/", "", "", "https://discord.heatlabs.net")

# Hyperlinks.Add() resets the cell's font to Excel's generic "Hyperlink"
# look; pin it back to match the sheet's established link style (row 249,
# etc.) so the style table doesn't grow with a near-duplicate entry.
$ws.Range("A250").Font.Name = "Arial"
$ws.Range("A250").Font.Size = 10
$ws.Range("A250").Font.Underline = $true
$ws.Range("A250").Font.Color = 16711680
